# Insert a new weekly price record as row 227 on the "Albahaca" (basil)
# price sheet, pushing the existing rows 227-283 down to 228-284.
#
# This mirrors a normal "insert row" edit made directly in Excel: the row
# is inserted (which shifts all subsequent rows down, carries formatting
# along, and grows the used range/dimension automatically), and then the
# new row is populated with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227; Excel shifts rows 227:283 down to
# 228:284 and copies the formatting (incl. the date style on column D)
# from the row above, same as an interactive "Insert" in the UI.
$ws.Rows.Item(227).EntireRow.Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A227").Value2 = 10
$ws.Range("B227").Value2 = "Vega Modelo de Temuco"
$ws.Range("C227").Value2 = "La Araucanía"
$ws.Range("D227").Value2 = 44889
$ws.Range("E227").Value2 = 9
$ws.Range("F227").Value2 = 100112052
$ws.Range("G227").Value2 = "Albahaca"
$ws.Range("H227").Value2 = "Sin especificar"
$ws.Range("I227").Value2 = "Primera"
$ws.Range("J227").Value2 = 55
$ws.Range("K227").Value2 = 9000
$ws.Range("L227").Value2 = 9000
$ws.Range("M227").Value2 = 9000
$ws.Range("N227").Value2 = "`$/paquete"
$ws.Range("O227").Value2 = "Región Metropolitana"
$ws.Range("P227").Value2 = 9000
$ws.Range("Q227").Value2 = 1
$ws.Range("R227").Value2 = "Hortaliza"
